$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2057.75
$ws.Range("I43").Value = 3160.2
$ws.Range("J43").Value = 1270.2858
$ws.Range("K43").Value = 3160.2
$ws.Range("L43").Value = 1270.2858
$ws.Range("M43").Value = -3091.2
$ws.Range("N43").Value = -1408.2858
$ws.Range("H61").Value = 69.166664
$ws.Range("I61").Value = 69.166664
$ws.Range("K61").Value = 207.499992
$ws.Range("M61").Value = -35.49999199999999
$ws.Range("H98").Value = 524.17145
$ws.Range("I98").Value = 524.17145
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 524.17145
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 4004.2432
$ws.Range("I100").Value = 2216.5625
$ws.Range("J100").Value = 5366.2856
$ws.Range("K100").Value = 2216.5625
$ws.Range("L100").Value = 5366.2856
$ws.Range("M100").Value = -1675.5625
$ws.Range("N100").Value = -6448.2856
$ws.Range("H122").Value = 524.17145
$ws.Range("I122").Value = 524.17145
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1572.51435
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4766.2163
$ws.Range("I132").Value = 2678.4
$ws.Range("J132").Value = 9115.833000000001
$ws.Range("K132").Value = 8035.200000000001
$ws.Range("L132").Value = 27347.499
$ws.Range("M132").Value = -5505.200000000001
$ws.Range("N132").Value = -32407.499

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3815.8857
$ws.Range("I61").Value = 2974
$ws.Range("J61").Value = 4707.294
$ws.Range("K61").Value = 2974
$ws.Range("L61").Value = 4707.294
$ws.Range("M61").Value = -2762
$ws.Range("N61").Value = -5131.294
$ws.Range("H97").Value = 660.5
$ws.Range("I97").Value = 669.94116
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 669.94116
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -173.94116
$ws.Range("N97").Value = -1492
$ws.Range("H102").Value = 2290.0667
$ws.Range("I102").Value = 1966
$ws.Range("J102").Value = 2938.2
$ws.Range("K102").Value = 1966
$ws.Range("L102").Value = 2938.2
$ws.Range("M102").Value = -344
$ws.Range("N102").Value = -6182.2
$ws.Range("H122").Value = 2031.4865
$ws.Range("I122").Value = 2391.318
$ws.Range("J122").Value = 1503.7333
$ws.Range("K122").Value = 7173.954000000001
$ws.Range("L122").Value = 4511.199900000001
$ws.Range("M122").Value = -4723.954000000001
$ws.Range("N122").Value = -9411.1999
$ws.Range("H132").Value = 3245.28
$ws.Range("I132").Value = 3229.4849
$ws.Range("J132").Value = 3275.9412
$ws.Range("K132").Value = 9688.4547
$ws.Range("L132").Value = 9827.8236
$ws.Range("M132").Value = -7158.4547
$ws.Range("N132").Value = -14887.8236
$ws.Range("H136").Value = 3815.8857
$ws.Range("I136").Value = 2974
$ws.Range("J136").Value = 4707.294
$ws.Range("K136").Value = 8922
$ws.Range("L136").Value = 14121.882
$ws.Range("M136").Value = -6372
$ws.Range("N136").Value = -19221.882

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1952.0555
$ws.Range("I20").Value = 1918.909
$ws.Range("J20").Value = 2004.1428
$ws.Range("K20").Value = 1918.909
$ws.Range("L20").Value = 2004.1428
$ws.Range("M20").Value = -1671.909
$ws.Range("N20").Value = -2498.1428
$ws.Range("H94").Value = 913.1667
$ws.Range("I94").Value = 716.3333
$ws.Range("J94").Value = 1110
$ws.Range("K94").Value = 716.3333
$ws.Range("L94").Value = 1110
$ws.Range("M94").Value = -265.3333
$ws.Range("N94").Value = -2012
$ws.Range("H96").Value = 13642.667
$ws.Range("I96").Value = 5464
$ws.Range("K96").Value = 5464
$ws.Range("M96").Value = -2718
$ws.Range("H107").Value = 2084.842
$ws.Range("I107").Value = 2243.5
$ws.Range("J107").Value = 1812.8572
$ws.Range("K107").Value = 2243.5
$ws.Range("L107").Value = 1812.8572
$ws.Range("M107").Value = -323.5
$ws.Range("N107").Value = -5652.8572
$ws.Range("H134").Value = 37904.5
$ws.Range("I134").Value = 68549
$ws.Range("J134").Value = 7260
$ws.Range("K134").Value = 205647
$ws.Range("L134").Value = 21780
$ws.Range("M134").Value = -203112
$ws.Range("N134").Value = -26850

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1950.5625
$ws.Range("I132").Value = 1138.9642
$ws.Range("J132").Value = 3086.8
$ws.Range("K132").Value = 3416.8926
$ws.Range("L132").Value = 9260.400000000001
$ws.Range("M132").Value = -886.8925999999997
$ws.Range("N132").Value = -14320.4
$ws.Range("H134").Value = 1823.2927
$ws.Range("I134").Value = 1255.3158
$ws.Range("J134").Value = 2313.818
$ws.Range("K134").Value = 3765.9474
$ws.Range("L134").Value = 6941.454000000001
$ws.Range("M134").Value = -1230.9474
$ws.Range("N134").Value = -12011.454

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 800
$ws.Range("I124").Value = 800
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 2400
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 4000
$ws.Range("I48").Value = 4000
$ws.Range("K48").Value = 4000
$ws.Range("M48").Value = -3515
$ws.Range("H102").Value = 1870.6471
$ws.Range("I102").Value = 2141.75
$ws.Range("J102").Value = 1220
$ws.Range("K102").Value = 2141.75
$ws.Range("L102").Value = 1220
$ws.Range("M102").Value = -519.75
$ws.Range("N102").Value = -4464
$ws.Range("H132").Value = 3712.319
$ws.Range("I132").Value = 3827.5518
$ws.Range("J132").Value = 3526.6667
$ws.Range("K132").Value = 11482.6554
$ws.Range("L132").Value = 10580.0001
$ws.Range("M132").Value = -8952.6554
$ws.Range("N132").Value = -15640.0001

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2530.7
$ws.Range("I7").Value = 2336.75
$ws.Range("J7").Value = 2660
$ws.Range("K7").Value = 2336.75
$ws.Range("L7").Value = 2660
$ws.Range("M7").Value = -2224.75
$ws.Range("N7").Value = -2884
$ws.Range("H16").Value = 1286.138
$ws.Range("I16").Value = 1279.0834
$ws.Range("J16").Value = 1320
$ws.Range("K16").Value = 1279.0834
$ws.Range("L16").Value = 1320
$ws.Range("M16").Value = -1109.0834
$ws.Range("N16").Value = -1660
$ws.Range("H22").Value = 401.81818
$ws.Range("I22").Value = 277.5
$ws.Range("K22").Value = 277.5
$ws.Range("M22").Value = 17.5
$ws.Range("H27").Value = 401.81818
$ws.Range("I27").Value = 277.5
$ws.Range("K27").Value = 277.5
$ws.Range("M27").Value = -170.5
$ws.Range("H40").Value = 2830.0588
$ws.Range("I40").Value = 2690.0908
$ws.Range("J40").Value = 3086.6667
$ws.Range("K40").Value = 2690.0908
$ws.Range("L40").Value = 3086.6667
$ws.Range("M40").Value = -2554.0908
$ws.Range("N40").Value = -3358.6667
$ws.Range("H46").Value = 704
$ws.Range("I46").Value = 404.66666
$ws.Range("J46").Value = 867.2727
$ws.Range("K46").Value = 404.66666
$ws.Range("L46").Value = 867.2727
$ws.Range("M46").Value = -216.66666
$ws.Range("N46").Value = -1243.2727
$ws.Range("H126").Value = 2530.7
$ws.Range("I126").Value = 2336.75
$ws.Range("J126").Value = 2660
$ws.Range("K126").Value = 7010.25
$ws.Range("L126").Value = 7980
$ws.Range("M126").Value = -4540.25
$ws.Range("N126").Value = -12920
$ws.Range("H128").Value = 20396
$ws.Range("J128").Value = 20396
$ws.Range("L128").Value = 20396
$ws.Range("N128").Value = -30356

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 51686.75
$ws.Range("I122").Value = 92525.91
$ws.Range("J122").Value = 1772.2222
$ws.Range("K122").Value = 277577.73
$ws.Range("L122").Value = 5316.6666
$ws.Range("M122").Value = -275127.73
$ws.Range("N122").Value = -10216.6666
$ws.Range("H136").Value = 36038976
$ws.Range("I136").Value = 52634188
$ws.Range("J136").Value = 18521808
$ws.Range("K136").Value = 157902564
$ws.Range("L136").Value = 55565424
$ws.Range("M136").Value = -157900014
$ws.Range("N136").Value = -55570524
